$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the batch/threshold related cells for socketTimeOut-related mysql test data
$ws.Range("I2").Value = "10000"
$ws.Range("L2").Value = "10000"
$ws.Range("M2").Value = "select count(*) from `$schema26 where id>2000 and id<=10000"
$ws.Range("N2").Value = "8000"
$ws.Range("O2").Value = "update `$schema26 set name='BJ' where id>2000 and id<=10000"
$ws.Range("P2").Value = "8000"
$ws.Range("R2").Value = "8000"
$ws.Range("T2").Value = "10000"

# Update the selected cell/active cell on the sheet
$ws.Activate()
$ws.Range("C10").Select()
